# Add a new worksheet named "ʤ" (U+02A4) at the end of the workbook,
# matching the csvkit fixture update that exercises utf-8 --sheet
# arguments. The sheet gets a simple 3-column table: headers a/b/c
# over one data row 1/2/3.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "ʤ"

$ws.Range("A1").Value = "a"
$ws.Range("B1").Value = "b"
$ws.Range("C1").Value = "c"
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 3

# Restore "data" as the active sheet/tab so activeTab stays 1, matching
# the target workbook (the new sheet is appended, not focused).
$null = $wb.Worksheets.Item("data").Activate()
